$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as text (matching the original inline-string
# storage) instead of being auto-parsed as numbers by the "50.72"-style values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.007.85"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.835.08"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "355.94"
$ws.Range("E5").Value = "  +7.06%  "
$ws.Range("D6").Value = "114.16"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").Value = "41.65"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "20.06"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "7.75"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "3.275.77"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "2.840.63"
$ws.Range("E16").Value = "  +3.96%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "51.931.11"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  +8.12%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "271.19"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "69.96"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +5.20%  "
$ws.Range("D26").Value = "26.80"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").Value = "0.0459"
$ws.Range("E31").Value = "  +30.75%  "
$ws.Range("D32").Value = "50.73"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "33.97"
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("E34").Value = "  +5.06%  "
$ws.Range("D35").Value = "0.0833"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "4.92"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "18.13"
$ws.Range("E40").Value = "  -4.90%  "
$ws.Range("D41").Value = "23.84"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("D44").Value = "126.59"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "2.092.56"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D49").Value = "5.73"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("D50").Value = "0.949"
$ws.Range("E50").Value = "  +8.88%  "
$ws.Range("D51").Value = "60.83"
$ws.Range("E51").Value = "  +1.55%  "

# Restore the default (General) cell style so no stray text-format override
# is left behind on the cells once the text values are in place.
$ws.Range("D2:D51").Style = "Normal"
